# Auto-generated: update price/profit columns (H-N) per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 3176327.5  # H17: 2470921.2 -> 3176327.5
$ws.Cells.Item(17, 10).Value = 3176327.5  # J17: 2470921.2 -> 3176327.5
$ws.Cells.Item(17, 12).Value = 9528982.5  # L17: 7412763.600000001 -> 9528982.5
$ws.Cells.Item(17, 14).Value = -9529318.5  # N17: -7413099.600000001 -> -9529318.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 4678.8335  # H74: 4734 -> 4678.8335
$ws.Cells.Item(74, 9).Value = 3599.375  # I74: 3626.4285 -> 3599.375
$ws.Cells.Item(74, 11).Value = 3599.375  # K74: 3626.4285 -> 3599.375
$ws.Cells.Item(74, 13).Value = -2663.375  # M74: -2690.4285 -> -2663.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 4678.8335  # H77: 4734 -> 4678.8335
$ws.Cells.Item(77, 9).Value = 3599.375  # I77: 3626.4285 -> 3599.375
$ws.Cells.Item(77, 11).Value = 17996.875  # K77: 18132.1425 -> 17996.875
$ws.Cells.Item(77, 13).Value = -13316.875  # M77: -13452.1425 -> -13316.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 887.25  # H98: 464.25925 -> 887.25
$ws.Cells.Item(98, 9).Value = 905.2  # I98: 437.6 -> 905.2
$ws.Cells.Item(98, 11).Value = 905.2  # K98: 437.6 -> 905.2
$ws.Cells.Item(98, 13).Value = 592.8  # M98: 1060.4 -> 592.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 887.25  # H122: 464.25925 -> 887.25
$ws.Cells.Item(122, 9).Value = 905.2  # I122: 437.6 -> 905.2
$ws.Cells.Item(122, 11).Value = 2715.6  # K122: 1312.8 -> 2715.6
$ws.Cells.Item(122, 13).Value = -265.6000000000004  # M122: 1137.2 -> -265.6000000000004

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 6087.8887  # H125: 6628.2856 -> 6087.8887
$ws.Cells.Item(125, 9).Value = 4196.5  # I125: 0 -> 4196.5
$ws.Cells.Item(125, 11).Value = 37768.5  # K125: 0 -> 37768.5
$ws.Cells.Item(125, 13).Value = -35308.5  # M125: None -> -35308.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2473.7307  # H137: 1817.1794 -> 2473.7307
$ws.Cells.Item(137, 9).Value = 2607.2942  # I137: 1695.9 -> 2607.2942
$ws.Cells.Item(137, 11).Value = 7821.882599999999  # K137: 5087.700000000001 -> 7821.882599999999
$ws.Cells.Item(137, 13).Value = -5271.882599999999  # M137: -2537.700000000001 -> -5271.882599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 0  # H10: 1100 -> 0
$ws.Cells.Item(10, 10).Value = 0  # J10: 1100 -> 0
$ws.Cells.Item(10, 12).Value = 0  # L10: 1100 -> 0
$ws.Cells.Item(10, 14).ClearContents()  # N10: -1440 -> (removed)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14092695  # H32: 14293862 -> 14092695
$ws.Cells.Item(32, 10).Value = 19912.863  # J32: 20336.715 -> 19912.863
$ws.Cells.Item(32, 12).Value = 19912.863  # L32: 20336.715 -> 19912.863
$ws.Cells.Item(32, 14).Value = -20486.863  # N32: -20910.715 -> -20486.863

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 11391.5  # H105: 12523.889 -> 11391.5
$ws.Cells.Item(105, 9).Value = 17690.334  # I105: 20988.4 -> 17690.334
$ws.Cells.Item(105, 11).Value = 17690.334  # K105: 20988.4 -> 17690.334
$ws.Cells.Item(105, 13).Value = -15943.334  # M105: -19241.4 -> -15943.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 19612728  # H31: 20004844 -> 19612728
$ws.Cells.Item(31, 10).Value = 52639110  # J31: 55563120 -> 52639110
$ws.Cells.Item(31, 12).Value = 52639110  # L31: 55563120 -> 52639110
$ws.Cells.Item(31, 14).Value = -52639700  # N31: -55563710 -> -52639700

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 19612728  # H34: 20004844 -> 19612728
$ws.Cells.Item(34, 10).Value = 52639110  # J34: 55563120 -> 52639110
$ws.Cells.Item(34, 12).Value = 52639110  # L34: 55563120 -> 52639110
$ws.Cells.Item(34, 14).Value = -52639514  # N34: -55563524 -> -52639514

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 58.53846  # H38: 53.357143 -> 58.53846
$ws.Cells.Item(38, 9).Value = 29.625  # I38: 26.777779 -> 29.625
$ws.Cells.Item(38, 10).Value = 104.8  # J38: 101.2 -> 104.8
$ws.Cells.Item(38, 11).Value = 88.875  # K38: 80.333337 -> 88.875
$ws.Cells.Item(38, 12).Value = 314.4  # L38: 303.6 -> 314.4
$ws.Cells.Item(38, 13).Value = 258.125  # M38: 266.666663 -> 258.125
$ws.Cells.Item(38, 14).Value = -1008.4  # N38: -997.6 -> -1008.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 20107  # H15: 69999 -> 20107
$ws.Cells.Item(15, 9).Value = 20107  # I15: 0 -> 20107
$ws.Cells.Item(15, 10).Value = 0  # J15: 69999 -> 0
$ws.Cells.Item(15, 11).Value = 20107  # K15: 0 -> 20107
$ws.Cells.Item(15, 12).Value = 0  # L15: 69999 -> 0
$ws.Cells.Item(15, 13).Value = -19819  # M15: None -> -19819
$ws.Cells.Item(15, 14).ClearContents()  # N15: -70575 -> (removed)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3058.4  # H80: 2984.2942 -> 3058.4
$ws.Cells.Item(80, 9).Value = 2769  # I80: 2779.111 -> 2769
$ws.Cells.Item(80, 10).Value = 3389.1428  # J80: 3215.125 -> 3389.1428
$ws.Cells.Item(80, 11).Value = 2769  # K80: 2779.111 -> 2769
$ws.Cells.Item(80, 12).Value = 3389.1428  # L80: 3215.125 -> 3389.1428
$ws.Cells.Item(80, 13).Value = -1771  # M80: -1781.111 -> -1771
$ws.Cells.Item(80, 14).Value = -5385.1428  # N80: -5211.125 -> -5385.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(81, 8).Value = 20107  # H81: 69999 -> 20107
$ws.Cells.Item(81, 9).Value = 20107  # I81: 0 -> 20107
$ws.Cells.Item(81, 10).Value = 0  # J81: 69999 -> 0
$ws.Cells.Item(81, 11).Value = 20107  # K81: 0 -> 20107
$ws.Cells.Item(81, 12).Value = 0  # L81: 69999 -> 0
$ws.Cells.Item(81, 13).Value = -19109  # M81: None -> -19109
$ws.Cells.Item(81, 14).ClearContents()  # N81: -71995 -> (removed)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3058.4  # H83: 2984.2942 -> 3058.4
$ws.Cells.Item(83, 9).Value = 2769  # I83: 2779.111 -> 2769
$ws.Cells.Item(83, 10).Value = 3389.1428  # J83: 3215.125 -> 3389.1428
$ws.Cells.Item(83, 11).Value = 13845  # K83: 13895.555 -> 13845
$ws.Cells.Item(83, 12).Value = 16945.714  # L83: 16075.625 -> 16945.714
$ws.Cells.Item(83, 13).Value = -8853  # M83: -8903.555 -> -8853
$ws.Cells.Item(83, 14).Value = -26929.714  # N83: -26059.625 -> -26929.714

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(84, 8).Value = 20107  # H84: 69999 -> 20107
$ws.Cells.Item(84, 9).Value = 20107  # I84: 0 -> 20107
$ws.Cells.Item(84, 10).Value = 0  # J84: 69999 -> 0
$ws.Cells.Item(84, 11).Value = 60321  # K84: 0 -> 60321
$ws.Cells.Item(84, 12).Value = 0  # L84: 209997 -> 0
$ws.Cells.Item(84, 13).Value = -55329  # M84: None -> -55329
$ws.Cells.Item(84, 14).ClearContents()  # N84: -219981 -> (removed)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(86, 8).Value = 47400  # H86: 49800 -> 47400
$ws.Cells.Item(86, 10).Value = 47400  # J86: 49800 -> 47400
$ws.Cells.Item(86, 12).Value = 47400  # L86: 49800 -> 47400
$ws.Cells.Item(86, 14).Value = -49772  # N86: -52172 -> -49772

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(89, 8).Value = 47400  # H89: 49800 -> 47400
$ws.Cells.Item(89, 10).Value = 47400  # J89: 49800 -> 47400
$ws.Cells.Item(89, 12).Value = 142200  # L89: 149400 -> 142200
$ws.Cells.Item(89, 14).Value = -154056  # N89: -161256 -> -154056

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2762.7222  # H132: 2834.7058 -> 2762.7222
$ws.Cells.Item(132, 9).Value = 2968.2  # I132: 3070.2856 -> 2968.2
$ws.Cells.Item(132, 11).Value = 8904.599999999999  # K132: 9210.856800000001 -> 8904.599999999999
$ws.Cells.Item(132, 13).Value = -6374.599999999999  # M132: -6680.856800000001 -> -6374.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3720.5557  # H7: 3521.718 -> 3720.5557
$ws.Cells.Item(7, 9).Value = 3124.7778  # I7: 2840.6191 -> 3124.7778
$ws.Cells.Item(7, 11).Value = 3124.7778  # K7: 2840.6191 -> 3124.7778
$ws.Cells.Item(7, 13).Value = -3012.7778  # M7: -2728.6191 -> -3012.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 2043.2  # H20: 2003.0385 -> 2043.2
$ws.Cells.Item(20, 9).Value = 1004.2222  # I20: 1003.7 -> 1004.2222
$ws.Cells.Item(20, 11).Value = 1004.2222  # K20: 1003.7 -> 1004.2222
$ws.Cells.Item(20, 13).Value = -778.2222  # M20: -777.7 -> -778.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(21, 8).Value = 3210  # H21: 5000 -> 3210
$ws.Cells.Item(21, 9).Value = 10000  # I21: 0 -> 10000
$ws.Cells.Item(21, 10).Value = 1512.5  # J21: 5000 -> 1512.5
$ws.Cells.Item(21, 11).Value = 10000  # K21: 0 -> 10000
$ws.Cells.Item(21, 12).Value = 1512.5  # L21: 5000 -> 1512.5
$ws.Cells.Item(21, 13).Value = -9826  # M21: None -> -9826
$ws.Cells.Item(21, 14).Value = -1860.5  # N21: -5348 -> -1860.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(24, 8).Value = 9666.666999999999  # H24: 10000 -> 9666.666999999999
$ws.Cells.Item(24, 9).Value = 12000  # I24: 0 -> 12000
$ws.Cells.Item(24, 10).Value = 5000  # J24: 10000 -> 5000
$ws.Cells.Item(24, 11).Value = 12000  # K24: 0 -> 12000
$ws.Cells.Item(24, 12).Value = 5000  # L24: 10000 -> 5000
$ws.Cells.Item(24, 13).Value = -11657  # M24: None -> -11657
$ws.Cells.Item(24, 14).Value = -5686  # N24: -10686 -> -5686

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2419.76  # H68: 2478.9167 -> 2419.76
$ws.Cells.Item(68, 9).Value = 2312.25  # I68: 2369.3044 -> 2312.25
$ws.Cells.Item(68, 11).Value = 2312.25  # K68: 2369.3044 -> 2312.25
$ws.Cells.Item(68, 13).Value = -1563.25  # M68: -1620.3044 -> -1563.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 2419.76  # H71: 2478.9167 -> 2419.76
$ws.Cells.Item(71, 9).Value = 2312.25  # I71: 2369.3044 -> 2312.25
$ws.Cells.Item(71, 11).Value = 11561.25  # K71: 11846.522 -> 11561.25
$ws.Cells.Item(71, 13).Value = -7817.25  # M71: -8102.522000000001 -> -7817.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(81, 8).Value = 92600  # H81: 83833.336 -> 92600
$ws.Cells.Item(81, 10).Value = 95000  # J81: 81250 -> 95000
$ws.Cells.Item(81, 12).Value = 95000  # L81: 81250 -> 95000
$ws.Cells.Item(81, 14).Value = -96996  # N81: -83246 -> -96996

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(84, 8).Value = 92600  # H84: 83833.336 -> 92600
$ws.Cells.Item(84, 10).Value = 95000  # J84: 81250 -> 95000
$ws.Cells.Item(84, 12).Value = 285000  # L84: 243750 -> 285000
$ws.Cells.Item(84, 14).Value = -294984  # N84: -253734 -> -294984

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(107, 8).Value = 10007  # H107: 16997.5 -> 10007
$ws.Cells.Item(107, 9).Value = 10007  # I107: 16997.5 -> 10007
$ws.Cells.Item(107, 11).Value = 10007  # K107: 16997.5 -> 10007
$ws.Cells.Item(107, 13).Value = -8087  # M107: -15077.5 -> -8087

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4332.6553  # H122: 3832.0881 -> 4332.6553
$ws.Cells.Item(122, 9).Value = 3802.7896  # I122: 3315.88 -> 3802.7896
$ws.Cells.Item(122, 10).Value = 5339.4  # J122: 5266 -> 5339.4
$ws.Cells.Item(122, 11).Value = 11408.3688  # K122: 9947.639999999999 -> 11408.3688
$ws.Cells.Item(122, 12).Value = 16018.2  # L122: 15798 -> 16018.2
$ws.Cells.Item(122, 13).Value = -8958.3688  # M122: -7497.639999999999 -> -8958.3688
$ws.Cells.Item(122, 14).Value = -20918.2  # N122: -20698 -> -20918.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 3720.5557  # H126: 3521.718 -> 3720.5557
$ws.Cells.Item(126, 9).Value = 3124.7778  # I126: 2840.6191 -> 3124.7778
$ws.Cells.Item(126, 11).Value = 9374.3334  # K126: 8521.8573 -> 9374.3334
$ws.Cells.Item(126, 13).Value = -6904.3334  # M126: -6051.8573 -> -6904.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 26888.834  # H95: 27197.75 -> 26888.834
$ws.Cells.Item(95, 10).Value = 26888.834  # J95: 27197.75 -> 26888.834
$ws.Cells.Item(95, 12).Value = 26888.834  # L95: 27197.75 -> 26888.834
$ws.Cells.Item(95, 14).Value = -32380.834  # N95: -32689.75 -> -32380.834

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 47671910  # H122: 50055476 -> 47671910
$ws.Cells.Item(122, 9).Value = 50055308  # I122: 55616932 -> 50055308
$ws.Cells.Item(122, 10).Value = 3998  # J122: 2349 -> 3998
$ws.Cells.Item(122, 11).Value = 150165924  # K122: 166850796 -> 150165924
$ws.Cells.Item(122, 12).Value = 11994  # L122: 7047 -> 11994
$ws.Cells.Item(122, 13).Value = -150163474  # M122: -166848346 -> -150163474
$ws.Cells.Item(122, 14).Value = -16894  # N122: -11947 -> -16894
